# "Resolvendo o problema do relatorio"
#
# The report sheet (Plan1) had a stray/incorrect row for "Kaka" (row 2).
# Remove that whole row - this shifts the "Fernando" row up from row 3 to
# row 2, automatically keeping Fernando's own cell values/types/formatting
# (name, matricula "154" as text, and its date-only number format).
#
# After the shift, correct the birth date on the surviving row to the
# right value (2001-08-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()

$ws.Range("C2").Value = 37120
